# Implement Handsontable for business.
# Remove the "supplier" and "business" header columns (X2, Y2) from the
# template sheet, leaving the cells empty but keeping their formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("X2:Y2").ClearContents()
$ws.Range("Y2").Select()
